$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 0.9107244398953713
$ws.Range("I3").Value = 0.008067229122586471
$ws.Range("K3").Value = 122.5961538461538

$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = 18
$ws.Range("S3").Value = 46
$ws.Range("T3").Value = 81
$ws.Range("U3").Value = 100

$ws.Range("V3").Value = 8683
$ws.Range("W3").Value = 8671
$ws.Range("X3").Value = 8643
$ws.Range("Y3").Value = 8608
$ws.Range("Z3").Value = 8589

$ws.Range("AF3").Value = 0.999309
$ws.Range("AG3").Value = 0.997928
$ws.Range("AH3").Value = 0.994706
$ws.Range("AI3").Value = 0.9906779999999999
$ws.Range("AJ3").Value = 0.988491
